# Add version 1 everywhere
#
# 1. Insert a new "version list" worksheet (with a single text cell "1")
#    right after "Export as TSV".
# 2. Insert a new "version" column at the very front of "Export as TSV",
#    shifting every other column right by one, re-creating comments at
#    their new positions, and adding a list-validation for the new column.

$wb = $excel.ActiveWorkbook
$tsv = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------
# Step 1: capture the existing header comments (they do NOT travel with
# a column insert, so we must remove + re-create them ourselves).
# ---------------------------------------------------------------------
$existingComments = [ordered]@{}
foreach ($cell in $tsv.UsedRange.Rows.Item(1).Cells) {
    if ($cell.Comment -ne $null) {
        $existingComments[$cell.Column] = $cell.Comment.Text()
    }
}

foreach ($cell in $tsv.UsedRange.Rows.Item(1).Cells) {
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
}

# ---------------------------------------------------------------------
# Step 2: insert the new "version list" worksheet right after the
# "Export as TSV" sheet, containing the single allowed value "1".
# ---------------------------------------------------------------------
$versionList = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tsv)
$versionList.Name = "version list"
$versionList.Range("A1").NumberFormat = "@"
$versionList.Range("A1").Value = "1"
$versionList.Range("A1").Style = "Normal"

# ---------------------------------------------------------------------
# Step 3: insert the new "version" column at the front of the TSV sheet.
# ---------------------------------------------------------------------
$tsv.Columns.Item(1).Insert()

# Copy the header formatting (bold, centered, wrap text) from the
# column that used to be A (now B) onto the new A1 cell.
$tsv.Range("B1").Copy()
$tsv.Range("A1").PasteSpecial(-4122) | Out-Null
$tsv.Application.CutCopyMode = $false
$tsv.Range("A1").Value = "version"

# New comment for the version column.
$tsv.Range("A1").AddComment("Current version of metadata schema. Template provides the correct value.") | Out-Null

# Re-create the rest of the header comments, shifted one column to the right.
foreach ($col in $existingComments.Keys) {
    $newCol = $col + 1
    $targetCell = $tsv.Cells.Item(1, $newCol)
    $targetCell.AddComment($existingComments[$col]) | Out-Null
}

# Data validation for the new "version" column.
$versionRange = $tsv.Range("A2:A1048576")
$versionRange.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1") | Out-Null
$versionRange.Validation.ErrorTitle = "Value must come from list"
$versionRange.Validation.ErrorMessage = "Value must be one of: 1."
$versionRange.Validation.ShowInput = $true
$versionRange.Validation.ShowError = $true

$tsv.Activate()
$tsv.Range("A1").Select() | Out-Null
